# Remove the 'PDB molecule' column from the metadata input file.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole "PDB molecule" column (column D); data to the right
# (PDB filename, etc.) shifts left automatically.
$ws.Range("D1:D8").EntireColumn.Delete()

# Widen the "Object name" column (now column B) a bit, as in the authored
# workbook (closest width this host can reproduce to 25.6640625 chars).
$ws.Columns("B").ColumnWidth = 24.8

# Leave the cell selection where Excel would land after the edit.
$ws.Range("G9").Select()
